# The edit re-orders the species-occurrence records that live in rows 4-15
# of the "Artfynd" sheet. Each row holds one full record (Id, TaxonId,
# species names, author, place name, coordinates, ...); the values that
# vary between records live in columns A,B,D,E,F,G,H,P,Q,R (C and the
# other columns are constant boiler-plate shared by every record). The
# edit moves whole records to new row positions (a permutation of rows
# 4-15) without touching rows 1-3 or anything outside columns A:S.
#
# We snapshot the "before" values for the columns that vary per record,
# then write them back out according to the row permutation below
# (destination row -> source/original row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry per-record data: A,B,D,E,F,G,H,P,Q,R
$cols = @(1, 2, 4, 5, 6, 7, 8, 16, 17, 18)

# Snapshot the current ("before") values for rows 4-15 across those columns,
# so later writes in the loop below never read already-overwritten data.
$snapshot = @{}
for ($r = 4; $r -le 15; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Row permutation: destination row -> source (original) row.
$mapping = @{}
$mapping[4]  = 5
$mapping[5]  = 7
$mapping[6]  = 8
$mapping[7]  = 9
$mapping[8]  = 10
$mapping[9]  = 11
$mapping[10] = 13
$mapping[11] = 14
$mapping[12] = 15
$mapping[13] = 4
$mapping[14] = 6
$mapping[15] = 12

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
